# Daily attendance processing - reorder the "Recorded By" (column G) list
# for each session row so it reflects the latest recorder-merge order.
# Rule observed in the day's export: the comma-separated list of recorders
# is reversed, except when it already ends with the literal entry "System"
# (single-entry cells are unaffected, since reversing them is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ", "
    if ($parts.Count -gt 1 -and $parts[$parts.Count - 1] -ne "System") {
        $reversed = $parts[($parts.Count - 1)..0]
        $cell.Value = ($reversed -join ", ")
    }
}
